$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlShiftDown = [Microsoft.Office.Interop.Excel.XlInsertShiftDirection]::xlShiftDown
$xlPasteFormats = [Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats

# Insert a new row 28 (pushes old blank rows / signature block down by one,
# table grows from 12 data rows [16-27] to 13 data rows [16-28])
$ws.Rows.Item(28).Insert($xlShiftDown)

# Copy the "closing/bottom border" formatting that used to sit on row 27
# onto the new row 28, so the border outline of the table moves to the new
# last row.
$ws.Range("B27:J27").Copy()
$ws.Range("B28:J28").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

# Re-populate the data table (now 13 rows, reordered newest period first,
# plus a new trailing row for the employer's own NIT/period 1606).
$data = @(
  @("CC","19391057","FERNANDO OROZCO ZAMORA","1903",26041,781242),
  @("CC","19391057","FERNANDO OROZCO ZAMORA","1902",31249,781242),
  @("CC","19391057","FERNANDO OROZCO ZAMORA","1901",31249,781242),
  @("CC","19391057","FERNANDO OROZCO ZAMORA","1812",31249,781242),
  @("CC","19391057","FERNANDO OROZCO ZAMORA","1811",31249,781242),
  @("CC","19391057","FERNANDO OROZCO ZAMORA","1612",27578,781242),
  @("CC","19391057","FERNANDO OROZCO ZAMORA","1611",27578,781242),
  @("CC","19391057","FERNANDO OROZCO ZAMORA","1610",27578,781242),
  @("CC","19391057","FERNANDO OROZCO ZAMORA","1609",27578,781242),
  @("CC","19391057","FERNANDO OROZCO ZAMORA","1608",27578,781242),
  @("CC","19391057","FERNANDO OROZCO ZAMORA","1607",27578,781242),
  @("CC","9089643","ANIANO ZABALA GOMEZ","1608",32000,800000)
)

$r = 16
foreach ($row in $data) {
  $ws.Cells.Item($r,2).Value = $row[0]
  $ws.Cells.Item($r,3).Value = $row[1]
  $ws.Cells.Item($r,4).Value = $row[2]
  $ws.Cells.Item($r,5).Value = $row[3]
  $ws.Cells.Item($r,6).Value = $row[4]
  $ws.Cells.Item($r,7).Value = $row[5]
  $r++
}

# New trailing row: the employer's own NIT, period 1606
$ws.Cells.Item(28,2).Value = "NIT"
$ws.Cells.Item(28,3).Value = "9009156792"
$ws.Cells.Item(28,4).Value = ""
$ws.Cells.Item(28,5).Value = "1606"
$ws.Cells.Item(28,6).Value = 234720
$ws.Cells.Item(28,7).Value = 0

# Update totals: total mora value (E11) and counts (C13 workers, F13 periods)
$ws.Cells.Item(11,5).Value = 583225
$ws.Cells.Item(13,3).Value = 3
$ws.Cells.Item(13,6).Value = 12

# Autofit the data columns whose widths changed as a side effect of the new content
for ($c = 2; $c -le 10; $c++) {
  $ws.Columns.Item($c).AutoFit()
}
